$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 7 - this shifts existing rows 7..131 down to 8..132
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new record's values.
$ws.Cells.Item(7, 1).Value = 7                                          # A7 - Mercado ID
$ws.Cells.Item(7, 2).Value = "Terminal Hortofrutícola Agro Chillán"     # B7 - Mercado
$ws.Cells.Item(7, 3).Value = "Ñuble"                                    # C7 - Región
$ws.Cells.Item(7, 4).Value = 45245                                      # D7 - Fecha
$ws.Cells.Item(7, 5).Value = 16                                         # E7 - Codreg
$ws.Cells.Item(7, 6).Value = 100112001                                  # F7 - Categoría ID
$ws.Cells.Item(7, 7).Value = "Berenjena"                                # G7 - Categoría
$ws.Cells.Item(7, 8).Value = "Sin especificar"                          # H7 - Variedad
$ws.Cells.Item(7, 9).Value = "Primera"                                  # I7 - Calidad
$ws.Cells.Item(7, 10).Value = 100                                       # J7 - Volumen
$ws.Cells.Item(7, 11).Value = 12000                                     # K7 - Precio mínimo
$ws.Cells.Item(7, 12).Value = 12000                                     # L7 - Precio máximo
$ws.Cells.Item(7, 13).Value = 12000                                     # M7 - Precio promedio ponderado
$ws.Cells.Item(7, 14).Value = "$/caja 60 unidades"                      # N7 - Unidad de comercialización
$ws.Cells.Item(7, 15).Value = "Región de Arica y Parinacota"            # O7 - Origen
$ws.Cells.Item(7, 16).Value = 200                                       # P7 - Precio $/Kg
$ws.Cells.Item(7, 17).Value = 60                                        # Q7 - Kg o Unidades
$ws.Cells.Item(7, 18).Value = "Hortaliza"                               # R7 - Clasificación
